$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab stays same sheetId/order, only the display name changes)
$ws.Name = "cases"

# Update existing row values (A3, A4)
$ws.Range("A3").Value2 = 1
$ws.Range("A4").Value2 = 2

# Add new row 5 data (new shared strings must be introduced in the same order
# they appear in the target sharedStrings table: dd_4, uieyfhr, then sdfg, awerdfs)
$ws.Range("A5").Value2 = 3
$ws.Range("C5").Value2 = "dd_4"
$ws.Range("D5").Value2 = "uieyfhr"
$ws.Range("E5").Value2 = 26
$ws.Range("G5").Value2 = "sam"

# H5 needs to be a date (43864) formatted like H3 (reuse same number format/style)
$ws.Range("H3").Copy() | Out-Null
$ws.Range("H5").PasteSpecial(-4122) | Out-Null
$ws.Range("H5").Value2 = 43864

$ws.Range("I5").Value2 = "sam"

# Update existing G4/I4 text values to new strings
$ws.Range("G4").Value2 = "sdfg"
$ws.Range("I4").Value2 = "awerdfs"

# Widen column G slightly (best achievable value close to the authored 16.5546875;
# the Excel object model quantizes ColumnWidth to whole pixels, so an exact match
# to the fractional authored width isn't reachable -- this lands on the nearest one)
$ws.Columns(7).ColumnWidth = 15.65

# Update the view's active cell/selection
$ws.Range("G9").Select() | Out-Null
